# feat/ added home page content and more recent rounds
#
# The underlying table (Summary sheet, A1:I.. ) is a golf-round log.
# This edit:
#   1) Adds a new "Trentham Golf Course" outing (5 players, row 31-35)
#   2) Re-sorts/re-applies the autofilter over A1:I35 by Date (column A)
#      ascending - this is what shuffles the existing rows 24-30 into
#      their new positions (a round that was appended out-of-order at
#      the bottom - old row 30, 30 Sep 2025 - sorts back up to row 24).
#   3) Updates the (hidden) _FilterDatabase defined name to match.
#   4) Appends one more, most-recent round (row 36) after the filtered
#      table, not yet part of the sort/filter range.
#   5) Widens column I (Comment) to fit the new comment text, and
#      updates the selection to where the user ended up (F38).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) New Trentham Golf Course rows (appended before the re-sort)
# ---------------------------------------------------------------------

$ws.Cells.Item(31,1).Value = 46046
$ws.Cells.Item(31,1).NumberFormat = "dd/mm/yy"
$ws.Cells.Item(31,2).Value = "Trentham Golf Course"
$ws.Cells.Item(31,3).Value = "Russell"
$ws.Cells.Item(31,4).Value = "Full-18"
$ws.Cells.Item(31,5).Value = "Solo"
$ws.Cells.Item(31,6).Value = 94
$ws.Cells.Item(31,7).Value = 71
$ws.Cells.Item(31,8).Formula = "=SUM(F31-G31)"

$ws.Cells.Item(32,1).Value = 46046
$ws.Cells.Item(32,1).NumberFormat = "dd/mm/yy"
$ws.Cells.Item(32,2).Value = "Trentham Golf Course"
$ws.Cells.Item(32,3).Value = "Olivia"
$ws.Cells.Item(32,4).Value = "Full-18"
$ws.Cells.Item(32,5).Value = "Solo"
$ws.Cells.Item(32,6).Value = 135
$ws.Cells.Item(32,7).Value = 71
$ws.Cells.Item(32,8).Formula = "=SUM(F32-G32)"

$ws.Cells.Item(33,1).Value = 46046
$ws.Cells.Item(33,1).NumberFormat = "dd/mm/yy"
$ws.Cells.Item(33,2).Value = "Trentham Golf Course"
$ws.Cells.Item(33,3).Value = "Hayden"
$ws.Cells.Item(33,4).Value = "Full-18"
$ws.Cells.Item(33,5).Value = "Solo"
$ws.Cells.Item(33,6).Value = 116
$ws.Cells.Item(33,7).Value = 71
$ws.Cells.Item(33,8).Formula = "=SUM(F33-G33)"

$ws.Cells.Item(34,1).Value = 46046
$ws.Cells.Item(34,1).NumberFormat = "dd/mm/yy"
$ws.Cells.Item(34,2).Value = "Trentham Golf Course"
$ws.Cells.Item(34,3).Value = "William"
$ws.Cells.Item(34,4).Value = "Full-18"
$ws.Cells.Item(34,5).Value = "Solo"
$ws.Cells.Item(34,6).Value = 147
$ws.Cells.Item(34,7).Value = 71
$ws.Cells.Item(34,8).Formula = "=SUM(F34-G34)"

$ws.Cells.Item(35,1).Value = 46046
$ws.Cells.Item(35,1).NumberFormat = "dd/mm/yy"
$ws.Cells.Item(35,2).Value = "Trentham Golf Course"
$ws.Cells.Item(35,3).Value = "Tohe"
$ws.Cells.Item(35,4).Value = "Full-18"
$ws.Cells.Item(35,5).Value = "Solo"
$ws.Cells.Item(35,6).Value = 89
$ws.Cells.Item(35,7).Value = 71
$ws.Cells.Item(35,8).Formula = "=SUM(F35-G35)"
$ws.Cells.Item(35,9).Value = "Bro’s so cracked"

# ---------------------------------------------------------------------
# 2) Re-sort A2:I35 by Date ascending + re-apply the autofilter
#    over the same A1:I35 range (matches the committed diff exactly -
#    this is what reorders the pre-existing rows 24-30 too).
# ---------------------------------------------------------------------

$ws.AutoFilterMode = $false
$ws.Range("A1:I35").AutoFilter()
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($ws.Range("A2:A35"))
$ws.Sort.SetRange($ws.Range("A1:I35"))
$ws.Sort.Header = 1
$ws.Sort.Apply()

# ---------------------------------------------------------------------
# 3) Keep the hidden _FilterDatabase name in sync with the new range
# ---------------------------------------------------------------------

foreach ($n in $wb.Names) {
    if ($n.Name -eq "Summary!_FilterDatabase") {
        $n.RefersTo = '=Summary!$A$1:$I$35'
    }
}

# ---------------------------------------------------------------------
# 4) Append the most recent round below the table (not yet part of
#    the sorted/filtered range - matches the source row's odd/default
#    styling on columns C:F).
# ---------------------------------------------------------------------

$ws.Cells.Item(36,1).Value = 46053
$ws.Cells.Item(36,1).NumberFormat = "dd/mm/yy"
$ws.Cells.Item(36,2).Value = "Mahunga"
$ws.Cells.Item(36,3).Value = "Golf Warehouse"
$ws.Cells.Item(36,4).Value = "Full-18"
$ws.Cells.Item(36,5).Value = "Solo"
$ws.Cells.Item(36,6).Formula = "=SUM(5+7+4+4+5+4+6+8+5+5+7+6+5+7+6+6+5+6)"
$ws.Cells.Item(36,7).Value = 71
$ws.Cells.Item(36,8).Formula = "=SUM(F36-G36)"
$ws.Cells.Item(36,9).Value = "This is the start of very honest golf. No Gimme’s, no Mulligans no breakfast balls or kick outs Wasn’t really trying here, too pissed off lol"

# ---------------------------------------------------------------------
# 5) Column width + selection / view
# ---------------------------------------------------------------------

$ws.Columns.Item(9).ColumnWidth = 113.85

$ws.Range("F38").Select()
